$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the journal entry date in D15 from 2021-03-11 to 2021-03-12
$ws.Range("D15").Value = 44267

# Move the active selection to D16, matching the saved selection state
$ws.Range("D16").Select()
